$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 161
$ws.Range("A2").Value = 153
$ws.Range("A3").Value = 177
$ws.Range("A4").Value = 147
$ws.Range("A5").Value = 135
$ws.Range("A6").Value = 175.6000000000004
$ws.Range("A7").Value = 161.3999999999996
$ws.Range("A8").Value = 194.2000000000007
